$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "'2025-10-16"
$ws.Range("B5").Value = "YYY"
$ws.Range("C5").Value = "123ABX007"
$ws.Range("D5").Value = "Karapakkam"
